$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 00:35"

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 1665984
$ws.Range("C4").Value = 20890
$ws.Range("D4").Value = 445377
$ws.Range("E4").Value = 1121951
$ws.Range("G4").Value = 1009
$ws.Range("H4").Value = 98656

# Colombia (row 39) - updated totals
$ws.Range("B39").Value = 20177
$ws.Range("C39").Value = 1046
$ws.Range("D39").Value = 4718
$ws.Range("E39").Value = 14754
$ws.Range("G39").Value = 23
$ws.Range("H39").Value = 705

# Argentina overtakes Dinamarca, Corea del Sur and Serbia (rows 47-50 shift down)
$ws.Range("A47").Value = "Argentina"
$ws.Range("B47").Value = 11353
$ws.Range("C47").Value = 704
$ws.Range("D47").Value = 3530
$ws.Range("E47").Value = 7378
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 445

$ws.Range("A48").Value = "Dinamarca"
$ws.Range("B48").Value = 11289
$ws.Range("C48").Value = 59
$ws.Range("D48").Value = 9836
$ws.Range("E48").Value = 892
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 561

$ws.Range("A49").Value = "Corea del Sur"
$ws.Range("B49").Value = 11165
$ws.Range("C49").Value = 23
$ws.Range("D49").Value = 10194
$ws.Range("E49").Value = 705
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 266

$ws.Range("A50").Value = "Serbia"
$ws.Range("B50").Value = 11092
$ws.Range("C50").Value = 68
$ws.Range("D50").Value = 5699
$ws.Range("E50").Value = 5155
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 238

# Maldivas (row 101) - updated active/recovered split
$ws.Range("D101").Value = 128
$ws.Range("E101").Value = 1181

# Mauritania overtakes Yemen and Birmania (rows 152-154 shift down)
$ws.Range("A152").Value = "Mauritania"
$ws.Range("B152").Value = 227
$ws.Range("C152").Value = 27
$ws.Range("D152").Value = 7
$ws.Range("E152").Value = 214
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 6

$ws.Range("A153").Value = "Yemen"
$ws.Range("B153").Value = 212
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 11
$ws.Range("E153").Value = 162
$ws.Range("G153").Value = 6
$ws.Range("H153").Value = 39

$ws.Range("A154").Value = "Birmania"
$ws.Range("B154").Value = 201
$ws.Range("C154").Value = 2
$ws.Range("D154").Value = 120
$ws.Range("E154").Value = 75
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 6
